$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update frame length / unit-converted values in column B
$ws.Range("B2").Value = 0.7
$ws.Range("B3").Value = 1
$ws.Range("B4").Value = 1.4
$ws.Range("B5").Value = 1.1000000000000001
$ws.Range("B6").Value = 2.2000000000000002
$ws.Range("B7").Value = 3.5

# Update the active selection to B5
$ws.Range("B5").Select()
